# Update task tracking - add entries for the "tet" task tracking days.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Attendent check")
$ws2 = $wb.Worksheets.Item("Deadline")

# ---- Sheet "Attendent check": fill in the 21/01/14 row (row 10) ----
# Attendance (black) / Late (yellow) / Attendance (black) / Attendance (black)
$ws1.Range("B6").Copy() | Out-Null
$ws1.Range("B10").PasteSpecial(-4122) | Out-Null
$ws1.Range("D6").Copy() | Out-Null
$ws1.Range("C10").PasteSpecial(-4122) | Out-Null
$ws1.Range("B6").Copy() | Out-Null
$ws1.Range("D10:E10").PasteSpecial(-4122) | Out-Null

# ---- Sheet "Attendent check": add a new tracking row for 23/01/04 (row 11) ----
$ws1.Range("A11").Value = "23/01/04"
$ws1.Range("B6").Copy() | Out-Null
$ws1.Range("B11:C11").PasteSpecial(-4122) | Out-Null
$ws1.Range("C6").Copy() | Out-Null
$ws1.Range("D11").PasteSpecial(-4122) | Out-Null
$ws1.Range("D6").Copy() | Out-Null
$ws1.Range("E11").PasteSpecial(-4122) | Out-Null

# ---- Sheet "Deadline": add a new tracking row for 23/01/14 (row 10) ----
$ws2.Range("A9").Copy() | Out-Null
$ws2.Range("A10").PasteSpecial(-4122) | Out-Null
$ws2.Range("A10").Value = "23/01/14"

$ws2.Range("B9:E9").Copy() | Out-Null
$ws2.Range("B10:E10").PasteSpecial(-4122) | Out-Null
$ws2.Range("B10:E10").Value = "Usecase"

$excel.CutCopyMode = 0

# ---- Update the active sheet / selections to match the edited state ----
$ws1.Range("E11").Select() | Out-Null
$ws2.Range("E10").Select() | Out-Null
$ws2.Activate() | Out-Null
